$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before the current row 11 (shifts existing rows 11.. down to 13..)
$ws.Rows("11:12").Insert()

# New row 11: general.bValueZeroThreshold
$ws.Range("A11").Value = "general.bValueZeroThreshold"
$ws.Range("D11").Value = "structural_preprocessing,reconstruction_diffusion"
$ws.Range("E11").Value = "numeric"
$ws.Range("F11").Value = "scalar nonempty nonnegative"
$ws.Range("G11").Value = "advanced"
$ws.Range("H11").Value = "B-values smaller or equal to this threshold are assumed to indicate b0-scans and set to b-value = 0."

# New row 12: general.bValueScalingTol
$ws.Range("A12").Value = "general.bValueScalingTol"
$ws.Range("D12").Value = "structural_preprocessing,reconstruction_diffusion"
$ws.Range("E12").Value = "numeric"
$ws.Range("F12").Value = "scalar nonempty nonnegative"
$ws.Range("G12").Value = "advanced"
$ws.Range("H12").Value = "B-vectors with a norm that deviates from 1 more than this threshold are labeled as potentially non-unit gradients."

# Apply the same style (F/G columns use style index 1 = numFmtId 49 "@" text format)
$ws.Range("F11:G12").Style = $ws.Range("F13").Style

# Column width updates for D and E
$ws.Columns("D").ColumnWidth = 15.1640625
$ws.Columns("E").ColumnWidth = 29.6640625

# View / window settings
$excel.ActiveWindow.WindowState = -4143
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H13").Select()
